# Applies the "Refined metadata to be additional tab" change:
#  1. Updates the F (time_taken) column timestamps on the "data" sheet.
#  2. Adds a new "metadata" worksheet (after "data") summarizing panel
#     query metadata, with a header row and a single data row.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Update time_taken (column F) values on the "data" sheet.
# ---------------------------------------------------------------------
$dataSheet.Range("F2").Value = "2021-10-05 14:35:47.947168"
$dataSheet.Range("F3").Value = "2021-10-05 14:35:47.947176"
$dataSheet.Range("F4").Value = "2021-10-05 14:35:47.947179"
$dataSheet.Range("F5").Value = "2021-10-05 14:35:47.947182"
$dataSheet.Range("F6").Value = "2021-10-05 14:35:47.947185"
$dataSheet.Range("F7").Value = "2021-10-05 14:35:47.947188"
$dataSheet.Range("F8").Value = "2021-10-05 14:35:47.947190"

# ---------------------------------------------------------------------
# 2. Add the "metadata" worksheet right after "data".
# ---------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Match the page margins used by the "data" sheet (values are in points;
# 0.75in/0.75in/1in/1in/0.5in/0.5in).
$metaSheet.PageSetup.LeftMargin = 54
$metaSheet.PageSetup.RightMargin = 54
$metaSheet.PageSetup.TopMargin = 72
$metaSheet.PageSetup.BottomMargin = 72
$metaSheet.PageSetup.HeaderMargin = 36
$metaSheet.PageSetup.FooterMargin = 36

# Reuse the header style (bold/border/centered) already present on the
# "data" sheet's header row, and the style used for the leading numeric
# index column, so no redundant style entries are created.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("E1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row.
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row.
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Susceptibility to Fungal Infections"
$metaSheet.Range("C2").Value = 236
$metaSheet.Range("E2").Value = "2021-04-06T21:09:09.864434Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:47.943388"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/236/?format=json"

# D2 ("1.0") must stay a text value, not be coerced into the number 1.
# Write it with a leading quote on a scratch cell (forcing text type),
# copy only the resulting value into D2, then discard the scratch row
# so no stray content remains on the sheet.
$scratch = $metaSheet.Range("Z100")
$scratch.Value = "'1.0"
$scratch.Copy()
$metaSheet.Range("D2").PasteSpecial(-4163)
$metaSheet.Rows.Item(100).Delete()
$excel.CutCopyMode = $false

# Keep "data" as the active/selected sheet, matching the original workbook
# (only the <sheets> list gained the new "metadata" entry).
$dataSheet.Activate() | Out-Null
